$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text, matching inlineStr cells in the source
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "257.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.54%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.55%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.566"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-5.54%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05898"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.19%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.629"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.90%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8547"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.38%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9363"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.11%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.01043"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1,610.01%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1386"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.11%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04873"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "35.55%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07073"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.80%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03066"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.60%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09113"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.37%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001523"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.53%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006028"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.17%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.492"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.17%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.189"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.07%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.64%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.72%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1270"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.77%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.911"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.02%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04273"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.90%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.25%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004288"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.10%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.04%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03822"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.73%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006250"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.67%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.44%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.05%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01386"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "32.20%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005371"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.15%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06589"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-39.58%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2523"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "11,768.39%"
